$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count + $used.Row - 1

$regex = [regex]'^(\d+) ч\. (\d+) мин\. (\d+) сек\.$'

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null) {
        $m = $regex.Match([string]$val)
        if ($m.Success) {
            $h = $m.Groups[1].Value
            $miNum = [int]$m.Groups[2].Value
            $seNum = [int]$m.Groups[3].Value
            $miStr = $miNum.ToString("D2")
            $seStr = $seNum.ToString("D2")
            $newVal = $h + " ч. " + $miStr + " мин. " + $seStr + " сек."
            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
